# JBSL-Twitter-Notes_Dirty.docx edit script
# 1) Justify (w:jc val="both") every existing paragraph.
# 2) Re-wrap the stray "mas" (in ". Cuando mas esecifico") with spell-check
#    proofErr markers, splitting that run into three runs.
# 3) Append the new "Optimización" section (dated 11/03/2024) at the end
#    of the document, each new paragraph justified like the rest.

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- 1) Justify every paragraph currently in the document -----------------
$wdAlignParagraphJustify = 3
$d.Paragraphs.Alignment = $wdAlignParagraphJustify

# --- 2) Fix the "Cuando mas esecifico" run, splitting out "mas" ----------
$lastParaIndex = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($lastParaIndex)
$targetRange = $targetPara.Range

$fixedParagraphXml = '<w:p xmlns:w="' + $wNs + '">' +
  '<w:pPr><w:jc w:val="both"/></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Cuanto mas complejo es un modelo, mas probabilidad hay de que ocurra </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>overfitting</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve">. Cuando </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>mas</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>esecifico</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> eres, puedes acertar </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>mas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> pero también puedes fallar mas</w:t></w:r>' +
  '</w:p>'

$targetRange.InsertXML($fixedParagraphXml)

# --- 3) Append the new "Optimización" section ------------------------------
function Add-ParagraphXml([string]$innerXml) {
    $endPos = $d.Content.End
    $insertionRange = $d.Range($endPos, $endPos)
    $paragraphXml = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:jc w:val="both"/></w:pPr>' + $innerXml + '</w:p>'
    $insertionRange.InsertXML($paragraphXml)
}

Add-ParagraphXml('<w:r><w:br w:type="page"/></w:r>')

Add-ParagraphXml('<w:r><w:lastRenderedPageBreak/><w:t>11/03/2024</w:t></w:r>')

Add-ParagraphXml('<w:r><w:t>Optimización</w:t></w:r>')

Add-ParagraphXml(
  '<w:r><w:t>Todos los modelos tienen una gran cantidad de parámetros para ajustar su comportamiento</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">. Probar todas las posibles combinaciones de modelos y parámetros puede resultar una tarea imposible. Para dar saltos muy grandes de rendimiento es necesario cambiar la metodología, sin </w:t></w:r>' +
  '<w:r><w:t>embargo</w:t></w:r>' +
  '<w:r><w:t>, para mejorar ligeramente el modelo, esta técnica resulta útil.</w:t></w:r>'
)

Add-ParagraphXml('<w:r><w:t>La optimización se suele realizar una vez satisfecho con el modelo y metodología ya seleccionado.</w:t></w:r>')

Add-ParagraphXml(
  '<w:r><w:t xml:space="preserve">Para realizar estas pruebas de optimización automáticamente, se puede utilizar la herramienta </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>GridSearchCV</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>. Entrena el modelo con todas las combinaciones y devuelve la mejor combinación.</w:t></w:r>'
)

Add-ParagraphXml(
  '<w:r><w:t xml:space="preserve">Esta metodología se puede combinar con la herramienta Pipeline para mejorar sus resultados. Es un </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>wrapper</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> para ejecutar varios pasos uno detrás de otros</w:t></w:r>'
)

Write-Output "edit complete"
